# Insert one new data row at row 9 (pushes the existing rows 9-38 down to
# 10-39) and populate it with a new price observation, matching the
# "Hortaliza, Feria Lagunitas de Puerto Montt - Alcachofa" weekly update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 9:38 down to 10:39, leaving row 9 free for the new record.
$ws.Rows("9:9").Insert()

# Fill in the new row 9 with the latest observation.
$ws.Range("A9").Value = 4
$ws.Range("B9").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C9").Value = "Los Lagos"
$ws.Range("D9").Value = 45107
$ws.Range("E9").Value = 10
$ws.Range("F9").Value = 100112013
$ws.Range("G9").Value = "Alcachofa"
$ws.Range("H9").Value = "Argentina(o)"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 60
$ws.Range("K9").Value = 20000
$ws.Range("L9").Value = 20000
$ws.Range("M9").Value = 20000
$ws.Range("N9").Value = '$/caja 50 unidades'
$ws.Range("O9").Value = "Provincia de Limarí"
$ws.Range("P9").Value = 400
$ws.Range("Q9").Value = 50
$ws.Range("R9").Value = "Hortaliza"
